$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.875.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.806.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.94"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.48%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3714"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07377"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8727"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.821.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.378"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.486"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07028"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008715"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  -3.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.871.07"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.61%  "
$ws.Range("E22").Value = "  -1.60%  "
$ws.Range("E23").Value = "  -2.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.020.52"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.895"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.59%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E27").Value = "  -1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.158"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.49%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.280"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "115.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08926"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7576"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.156"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.453"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.908"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01959"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05259"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.927"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.248"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.60%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.372"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5283"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("E44").Value = "  -2.91%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.518"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.56%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4987"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.53%  "
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.666"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06296"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "
